# Adds rows 1501-1561 to Sheet1 (new "08/07/2024" / "09/07/2024" / "10/07/2024"
# FII/DII watchlist blocks) and extends the sheet dimension accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# These two new cells hold ambiguous dd/mm/yyyy-looking text ("09/07/2024",
# "10/07/2024"). Pre-formatting as Text keeps them literal strings instead of
# being auto-converted to date serials, matching the existing date rows above
# (e.g. A1500 = "08/07/2024") already stored as plain text in this sheet.
$ws.Range("A1527").NumberFormat = "@"
$ws.Range("A1561").NumberFormat = "@"

# Build the A1501:J1561 block in one shot via a 2-D array write.
$data = New-Object "object[,]" 61,10
$data[0,0] = 'Buying Opportunity'
$data[0,1] = 'support Zone'
$data[0,2] = 'long buildup'
$data[0,3] = 'Short buildup'
$data[0,4] = 'FII ENTERING'
$data[1,0] = 'ALKEM'
$data[1,1] = 'AEROFLEX'
$data[1,2] = 'COLPAL'
$data[1,3] = 'IDFCFIRSTB'
$data[1,4] = 'ALKEM'
$data[1,5] = 5324.15
$data[1,6] = 156.75
$data[1,7] = 2995.55
$data[1,8] = 79.19
$data[1,9] = 5324.15
$data[2,0] = 'CESC'
$data[2,1] = 'ARVINDFASN'
$data[2,2] = 'ITC'
$data[2,3] = 'JINDALSTEL'
$data[2,4] = 'COLPAL'
$data[2,5] = 181.94
$data[2,6] = 519.25
$data[2,7] = 452.6
$data[2,8] = 1021.85
$data[2,9] = 2995.55
$data[3,0] = 'COLPAL'
$data[3,1] = 'BIRLACABLE'
$data[3,4] = 'ITC'
$data[3,5] = 2995.55
$data[3,6] = 246.22
$data[3,9] = 452.6
$data[4,0] = 'CONSUMBEES'
$data[4,1] = 'BIRLACORPN'
$data[4,4] = 'LT'
$data[4,5] = 125.32
$data[4,6] = 1575
$data[4,9] = 3666.1
$data[5,0] = 'DEEPAKNTR'
$data[5,1] = 'BIRLAMONEY'
$data[5,4] = 'NAUKRI'
$data[5,5] = 2700.9
$data[5,6] = 161.31
$data[5,9] = 6944.6
$data[6,0] = 'GLAND'
$data[6,1] = 'DENORA'
$data[6,5] = 2031.8
$data[6,6] = 1844.65
$data[7,0] = 'GMRP&UI'
$data[7,1] = 'EXXARO'
$data[7,5] = 94.73999999999999
$data[7,6] = 91.09
$data[8,0] = 'GREENPOWER'
$data[8,1] = 'FORTIS'
$data[8,5] = 22.54
$data[8,6] = 455.3
$data[9,0] = 'HCC'
$data[9,1] = 'GALAXYSURF'
$data[9,5] = 51.14
$data[9,6] = 2988.95
$data[10,0] = 'HINDCOMPOS'
$data[10,1] = 'HONASA'
$data[10,5] = 604.9
$data[10,6] = 455.7
$data[11,0] = 'ITC'
$data[11,1] = 'JINDALPOLY'
$data[11,5] = 452.6
$data[11,6] = 790.15
$data[12,0] = 'KANORICHEM'
$data[12,1] = 'KPRMILL'
$data[12,5] = 136.12
$data[12,6] = 867.1
$data[13,0] = 'KCPSUGIND'
$data[13,1] = 'MANUGRAPH'
$data[13,5] = 46.86
$data[13,6] = 21.91
$data[14,0] = 'KELLTONTEC'
$data[14,1] = 'MOTILALOFS'
$data[14,5] = 141.46
$data[14,6] = 543
$data[15,0] = 'KSCL'
$data[15,1] = 'MUFIN'
$data[15,5] = 993.05
$data[15,6] = 117.48
$data[16,0] = 'LUXIND'
$data[16,1] = 'NAGREEKEXP'
$data[16,5] = 1687.85
$data[16,6] = 36.88
$data[17,0] = 'MADRASFERT'
$data[17,1] = 'NCLIND'
$data[17,5] = 130.08
$data[17,6] = 245.19
$data[18,0] = 'MANAKALUCO'
$data[18,1] = 'NECCLTD'
$data[18,5] = 32.35
$data[18,6] = 28
$data[19,0] = 'MRPL'
$data[19,1] = 'ORISSAMINE'
$data[19,5] = 232.82
$data[19,6] = 7205
$data[20,0] = 'NAUKRI'
$data[20,1] = 'ROTO'
$data[20,5] = 6944.6
$data[20,6] = 651.1
$data[21,0] = 'PASUPTAC'
$data[21,5] = 40.94
$data[22,0] = 'PFS'
$data[22,5] = 51.51
$data[23,0] = 'PHOENIXLTD'
$data[23,5] = 4025.65
$data[24,0] = 'RUSTOMJEE'
$data[24,5] = 684.2
$data[25,0] = 'SADHNANIQ'
$data[25,5] = 76.81
$data[26,0] = '09/07/2024'
$data[27,0] = 'Buying Opportunity'
$data[27,1] = 'support Zone'
$data[27,2] = 'long buildup'
$data[27,3] = 'Short buildup'
$data[27,4] = 'FII ENTERING'
$data[28,0] = 'AXISBANK'
$data[28,1] = 'AARVI'
$data[28,3] = 'HAL'
$data[28,4] = 'BRITANNIA'
$data[28,5] = 1291.65
$data[28,6] = 135.06
$data[28,8] = 5486.15
$data[28,9] = 5755.55
$data[29,0] = 'BRITANNIA'
$data[29,1] = 'ADANIENT'
$data[29,5] = 5755.55
$data[29,6] = 3096
$data[30,0] = 'COLPAL'
$data[30,1] = 'AGI'
$data[30,5] = 3040.5
$data[30,6] = 854.65
$data[31,0] = 'CONCORDBIO'
$data[31,1] = 'ALKYLAMINE'
$data[31,5] = 1702.55
$data[31,6] = 2056.8
$data[32,0] = 'HEXATRADEX'
$data[32,1] = 'AURUM'
$data[32,5] = 187
$data[32,6] = 155.05
$data[33,0] = 'LUXIND'
$data[33,1] = 'BANSWRAS'
$data[33,5] = 1759.7
$data[33,6] = 167.12
$data[34,0] = 'MANAKALUCO'
$data[34,1] = 'BEPL'
$data[34,5] = 35.58
$data[34,6] = 147.58
$data[35,1] = 'CENTENKA'
$data[35,6] = 569.3
$data[36,1] = 'COCHINSHIP'
$data[36,6] = 2727.25
$data[37,1] = 'COFORGE'
$data[37,6] = 5669.1
$data[38,1] = 'DCAL'
$data[38,6] = 167.56
$data[39,1] = 'DREDGECORP'
$data[39,6] = 1304.45
$data[40,1] = 'DYNAMATECH'
$data[40,6] = 7017.75
$data[41,1] = 'EVERESTIND'
$data[41,6] = 1161.55
$data[42,1] = 'FAIRCHEMOR'
$data[42,6] = 1376.55
$data[43,1] = 'GALAXYSURF'
$data[43,6] = 2907.8
$data[44,1] = 'GRSE'
$data[44,6] = 2471.2
$data[45,1] = 'GRWRHITECH'
$data[45,6] = 2338.95
$data[46,1] = 'GTLINFRA'
$data[46,6] = 3.53
$data[47,1] = 'HINDZINC'
$data[47,6] = 659.2
$data[48,1] = 'IFBIND'
$data[48,6] = 1589.3
$data[49,1] = 'KAMATHOTEL'
$data[49,6] = 209.22
$data[50,1] = 'KBCGLOBAL'
$data[50,6] = 1.94
$data[51,1] = 'KRISHANA'
$data[51,6] = 284.5
$data[52,1] = 'MAZDOCK'
$data[52,6] = 5358.15
$data[53,1] = 'MIRZAINT'
$data[53,6] = 44.67
$data[54,1] = 'MTARTECH'
$data[54,6] = 1954.45
$data[55,1] = 'NRAIL'
$data[55,6] = 487.7
$data[56,1] = 'ONEPOINT'
$data[56,6] = 68.47
$data[57,1] = 'PRECOT'
$data[57,6] = 478.45
$data[58,1] = 'PRITIKAUTO'
$data[58,6] = 26.98
$data[59,1] = 'RKFORGE'
$data[59,6] = 905.85
$data[60,0] = '10/07/2024'

$ws.Range("A1501:J1561").Value = $data
